# Applies the cryptos-list refresh described in the commit:
# "Updated cryptos list on Wed May 17 09:53:31 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text that can look numeric ("1.000", "5.310", ...).
# Force text format first so Excel keeps the literal string (incl. trailing
# zeros / thousand-dot formatting) exactly like the original inline strings.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.865.92'
$ws.Range('E2').Value = '  -1.36%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.807.05'
$ws.Range('E3').Value = '  -1.01%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').Value = '  -0.43%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '310.23'
$ws.Range('E6').Value = '  -0.39%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4465'
$ws.Range('E7').Value = '  +5.06%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3673'
$ws.Range('E8').Value = '  -0.89%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07406'
$ws.Range('E9').Value = '  +2.15%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8556'
$ws.Range('E10').Value = '  -0.92%  '
$ws.Range('E11').Value = '  -1.81%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.804.13'
$ws.Range('E12').Value = '  -1.27%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.611'
$ws.Range('E13').Value = '  -1.76%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '92.54'
$ws.Range('E14').Value = '  +3.45%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.310'
$ws.Range('E15').Value = '  -0.22%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.07075'
$ws.Range('E16').Value = '  -0.17%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.001'
$ws.Range('E17').Value = '  -0.44%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008741'
$ws.Range('E18').Value = '  -1.43%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.000'
$ws.Range('E19').Value = '  -0.40%  '
$ws.Range('E20').Value = '  -1.40%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '26.894.36'
$ws.Range('E21').Value = '  -1.60%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.161'
$ws.Range('E22').Value = '  +0.41%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.85'
$ws.Range('E23').Value = '  -0.48%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.992'
$ws.Range('E24').Value = '  +0.17%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '151.85'
$ws.Range('E25').Value = '  -0.47%  '
$ws.Range('B26').Value = 'LidoDAOToken'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.187'
$ws.Range('E26').Value = '  -0.86%  '
$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.49'
$ws.Range('E27').Value = '  +0.40%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '5.213'
$ws.Range('E28').Value = '  -0.68%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '116.52'
$ws.Range('E29').Value = '  +0.00%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.08830'
$ws.Range('E30').Value = '  -0.26%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.7538'
$ws.Range('E31').Value = '  -0.52%  '
$ws.Range('E32').Value = '  -1.85%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.927'
$ws.Range('E33').Value = '  +5.12%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.460'
$ws.Range('E34').Value = '  -0.01%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.9997'
$ws.Range('E35').Value = '  -0.49%  '
$ws.Range('E36').Value = '  -2.52%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.01971'
$ws.Range('E37').Value = '  -0.23%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.05202'
$ws.Range('E38').Value = '  -1.29%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.5334'
$ws.Range('E39').Value = '  +5.37%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.870'
$ws.Range('E40').Value = '  +0.02%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '7.004'
$ws.Range('E41').Value = '  -4.70%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1691'
$ws.Range('E42').Value = '  -0.47%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.5179'
$ws.Range('E43').Value = '  +8.88%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.449'
$ws.Range('E44').Value = '  -3.26%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '10.55'
$ws.Range('E45').Value = '  -1.49%  '
$ws.Range('B46').Value = 'RenderToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.980'
$ws.Range('E46').Value = '  +6.39%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '105.46'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.674'
$ws.Range('E48').Value = '  +0.07%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.9993'
$ws.Range('E49').Value = '  -0.49%  '
$ws.Range('E50').Value = '  -0.86%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.9213'
$ws.Range('E51').Value = '  +0.65%  '
